# Update generator_info.xlsx: strip the "009-" prefix from meter serial
# numbers in column A, mark that column as Text so the values (which look
# numeric, e.g. 980E1F / 980E29 which Excel would otherwise read as
# scientific notation) are preserved exactly, widen the column, move the
# selection, and force portrait print orientation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (Meter Serial No.) becomes Text-formatted so alphanumeric serials
# display verbatim. Apply the format BEFORE writing the values so Excel
# doesn't try to interpret strings like "980E29" as scientific notation.
$ws.Columns.Item(1).NumberFormat = "@"

# New meter-serial values (formerly "009-XXXXX") for rows 3-10.
$ws.Range("A3").Value  = "980E1F"
$ws.Range("A4").Value  = "980DCD"
$ws.Range("A5").Value  = "980E2A"
$ws.Range("A6").Value  = "980E29"
$ws.Range("A7").Value  = "980B76"
$ws.Range("A8").Value  = "980B1E"
$ws.Range("A9").Value  = "980B1C"
$ws.Range("A10").Value = "980B13"

# Widen the column to fit the new style/content.
$ws.Columns.Item(1).ColumnWidth = 17.7213541666666

# Move the active selection to A11 (below the table).
[void]$ws.Range("A11").Select()

# Always print in portrait orientation.
$ws.PageSetup.Orientation = 1
